$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the source data, which lands at
# row 384 (pushing the existing rows 384-402 down to 385-403).
$ws.Rows.Item(384).Insert()

$ws.Range("A384").Value2 = 9
$ws.Range("B384").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C384").Value2 = "Metropolitana"
$ws.Range("D384").Value2 = 44706
$ws.Range("E384").Value2 = 13
$ws.Range("F384").Value2 = 100112039
$ws.Range("G384").Value2 = "Ciboulette"
$ws.Range("H384").Value2 = "Sin especificar"
$ws.Range("I384").Value2 = "Primera"
$ws.Range("J384").Value2 = 670
$ws.Range("K384").Value2 = 800
$ws.Range("L384").Value2 = 1000
$ws.Range("M384").Value2 = 896
$ws.Range("N384").Value2 = "`$/docena de atados"
$ws.Range("O384").Value2 = "Región Metropolitana"
$ws.Range("P384").Value2 = 299
$ws.Range("Q384").Value2 = 3
$ws.Range("R384").Value2 = "Hortaliza"
